$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet "Data" -> "Summary" ---------------------------
$ws.Name = "Summary"

# Re-assert formatting on the two pre-existing, untouched cells so a
# save/reload round-trip through the COM layer doesn't silently drop their
# named-style look (font size / bold) - "name" style (18pt) and "title"
# style (bold).
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# --- 2. Remove the old row 5 / row 6 content ----------------------------
# (Micro/SMEs/MSMEs header + "Enterprises (% of total)" row) — this content
# reappears further down the sheet at its new location, so just delete the
# two old rows outright (shifts nothing else, since they were the last rows).
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# --- 3. Create the new "title_" cell style (bold + underlined) ---------
$titleUnderline = $wb.Styles.Add("title_")
$titleUnderline.Font.Bold = $true
$titleUnderline.Font.Underline = $true

# --- 4. Populate the new layout -----------------------------------------

# Row 9 - new "Source Type" sub-heading, bold + underlined
$ws.Range("A9").Value = "Source Type: SME Associations"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# Row 11 - column headers (bold, same "title" look as before)
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Row 12 - data row
$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").Value = "'90"

# Row 13 - source note (italic, same "source" look as before)
$ws.Range("A13").Value = "Source: OBG, 2010"
$ws.Range("A13").Font.Italic = $true

# Row 19 - "OBG" heading (bold)
$ws.Range("A19").Value = "OBG"
$ws.Range("A19").Font.Bold = $true

# Row 20 - full citation (italic)
$ws.Range("A20").Value = 'Oxford Business Group (OBG), "Economic Update, Ghana boosts financial support for SME development by Oxford Business Group", 2014. Available at http://www.oxfordbusinessgroup.com/economic_updates/ghana-boosts-financial-support-sme-development'
$ws.Range("A20").Font.Italic = $true
